{"js": "// Apply the text edits described by the diff.\n// Paragraphs 0-5 get their text replaced in place; two brand-new\n// paragraphs are inserted right before the final (URL) paragraph,\n// whose own text is then updated to the new arXiv link.\n\nconst newTexts = [\n  \"\u26a1\ufe0f\ud83d\ude80\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05d5\u05de\u05d9 \u05e9\u05dc \u05de\u05d9\u05d9\u05e7 25.06.24:\u26a1\ufe0f\ud83d\ude80\",\n  \"Improving Reinforcement Learning from Human Feedback with Efficient Reward Model Ensemble\",\n  \"\u05d4\u05e1\u05e7\u05d9\u05e8\u05d4 \u05d4\u05d6\u05d5 \u05de\u05de\u05e9\u05d9\u05db\u05d4 \u05d0\u05ea \u05e7\u05d5 \u05d4\u05e1\u05e7\u05d9\u05e8\u05d5\u05ea \u05e2\u05dc \u05d4\u05de\u05d0\u05de\u05e8\u05d9\u05dd \u05e9\u05de\u05e0\u05e1\u05d9\u05dd \u05dc\u05e9\u05e4\u05e8 \u05e9\u05d9\u05d8\u05d5\u05ea RLHF \u05dc\u05d8\u05d9\u05d5\u05d1 (instruction tuning \u05d0\u05d5 \u05e4\u05e9\u05d5\u05d8 fine-tuning) \u05e9\u05dc \u05de\u05d5\u05d3\u05dc\u05d9 \u05e9\u05e4\u05d4. \u05d1\u05d7\u05dc\u05e7 \u05e9\u05dc \u05e9\u05d9\u05d8\u05ea RLHF (\u05dc\u05de\u05e9\u05dc PPO) \u05d0\u05e0\u05d5 \u05de\u05d0\u05de\u05e0\u05d9\u05dd \u05de\u05d5\u05d3\u05dc reward \u05de\u05d1\u05d5\u05e1\u05e1 \u05e2\u05dc \u05e1\u05d8 \u05e9\u05dc \u05e9\u05d0\u05dc\u05d5\u05ea \u05d5\u05ea\u05e9\u05d5\u05d1\u05d5\u05ea \u05de\u05d3\u05d5\u05e8\u05d2\u05d5\u05ea \u05e2\u05dc \u05d9\u05d3\u05d9 \u05d4\u05de\u05ea\u05d9\u05d9\u05d2\u05d9\u05dd \u05d4\u05d0\u05e0\u05d5\u05e9\u05d9\u05d9\u05dd. \u05de\u05d8\u05e8\u05d4 \u05e9\u05dc \u05de\u05d5\u05d3\u05dc \u05d6\u05d4 \u05dc\u05e1\u05e4\u05e7 \u05e6\u05d9\u05d5\u05df \u05dc\u05d6\u05d5\u05d2 (\u05e9\u05d0\u05dc\u05d4, \u05ea\u05e9\u05d5\u05d1\u05d4) \u05db\u05d0\u05e9\u05e8 \u05e6\u05d9\u05d5\u05df \u05d2\u05d1\u05d5\u05d4 \u05de\u05e6\u05d1\u05d9\u05e2 \u05e2\u05dc \u05ea\u05e9\u05d5\u05d1\u05d4 \u05d8\u05d5\u05d1\u05d4 \u05d5\u05e8\u05e6\u05d5\u05d9\u05d4. \u05dc\u05d0\u05d7\u05e8 \u05db\u05df \u05d0\u05e0\u05d5 \u05de\u05d0\u05de\u05e0\u05d9\u05dd (\u05de\u05d8\u05d9\u05d9\u05d1\u05d9\u05dd) \u05de\u05d5\u05d3\u05dc \u05e9\u05e4\u05d4 \u05db\u05d0\u05e9\u05e8 \u05d4\u05de\u05d8\u05e8\u05d4 \u05d4\u05d9\u05d0 \u05de\u05e7\u05e1\u05d5\u05dd \u05e9\u05dc \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d4 reward \u05ea\u05d5\u05da \u05e9\u05de\u05d9\u05e8\u05ea \u05e9\u05dc \u05de\u05e9\u05e7\u05dc\u05d9 \u05d4\u05de\u05d5\u05d3\u05dc \u05dc\u05de\u05e9\u05e7\u05dc\u05d9\u05dd \u05e9\u05d4\u05ea\u05d7\u05dc\u05e0\u05d5 \u05de\u05d4\u05dd (\u05e0\u05de\u05d3\u05d3 \u05e2\u05dc \u05d9\u05d3\u05d9 KL divergence \u05d1\u05d9\u05df \u05d4\u05ea\u05e4\u05dc\u05d2\u05d5\u05d9\u05d5\u05ea \u05d4\u05d8\u05d5\u05e7\u05e0\u05d9\u05dd \u05e9\u05dc \u05e9\u05e0\u05d9 \u05d4\u05de\u05d5\u05d3\u05dc\u05d9\u05dd). \u05db\u05dc \u05d6\u05d4 \u05de\u05ea\u05d1\u05e6\u05e2 on-the-fly \u05db\u05d0\u05e9\u05e8 \u05d4\u05d3\u05d5\u05d2\u05de\u05d0\u05d5\u05ea \u05e0\u05d5\u05e6\u05e8\u05d5\u05ea \u05e2\u05dc\u05d9 \u05d9\u05d3\u05d9 \u05d4\u05d2\u05e8\u05e1\u05d4 \u05d4\u05e2\u05d3\u05db\u05e0\u05d9\u05ea \u05e9\u05dc \u05d4\u05de\u05d5\u05d3\u05dc \u05d1\u05de\u05d4\u05dc\u05da \u05d4\u05d0\u05d9\u05de\u05d5\u05df.\",\n  \"\u05d4\u05d1\u05e2\u05d9\u05d4 \u05e2\u05dd \u05d4\u05d2\u05d9\u05e9\u05d4 \u05d4\u05d9\u05d0 reward hacking \u05db\u05d0\u05e9\u05e8 \u05dc\u05de\u05e8\u05d5\u05ea \u05d0\u05d9\u05d1\u05e8 \u05d4\u05e8\u05d2\u05d5\u05dc\u05e8\u05d9\u05d6\u05e6\u05d9\u05d4 (KL) \u05d4\u05de\u05d5\u05d3\u05dc \u05de\u05ea\u05db\u05e0\u05e1 \u05dc\u05de\u05e9\u05e7\u05dc\u05d9\u05dd \u05e9\u05de\u05d2\u05d9\u05e2\u05d9\u05dd \u05dc\u05e2\u05e8\u05db\u05d9\u05dd \u05d2\u05d1\u05d5\u05d4\u05d9\u05dd \u05e9\u05dc \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d9\u05ea reward \u05db\u05d0\u05e9\u05e8 \u05d4\u05de\u05d5\u05d3\u05dc \u05e2\u05e6\u05de\u05d5 \u05f4\u05dc\u05d0 \u05de\u05e1\u05e4\u05e7 \u05d0\u05ea \u05d4\u05e1\u05d7\u05d5\u05e8\u05d4\u05f4. \u05d4\u05de\u05d0\u05de\u05e8 \u05de\u05e6\u05d9\u05e2 \u05dc\u05d4\u05e9\u05ea\u05de\u05e9 \u05d1\u05db\u05de\u05d4 \u05de\u05d5\u05d3\u05dc\u05d9 reward \u05db\u05d9 ensemble \u05d6\u05d4 \u05ea\u05de\u05d9\u05d3 \u05d8\u05d5\u05d1. \u05d4\u05d1\u05e2\u05d9\u05d4 \u05e9\u05dc\u05d4\u05d7\u05d6\u05d9\u05e7 \u05d9\u05d5\u05ea\u05e8 \u05de\u05de\u05d5\u05d3\u05dc \u05d0\u05d7\u05d3 \u05d1\u05d6\u05de\u05df \u05d4\u05d0\u05d9\u05de\u05d5\u05df \u05d6\u05d4 \u05d9\u05e7\u05e8 \u05de\u05d1\u05d7\u05d9\u05e0\u05ea \u05d4\u05de\u05e9\u05d0\u05d1\u05d9\u05dd. \u05d4\u05de\u05d0\u05de\u05e8 \u05de\u05e6\u05d9\u05e2 \u05e9\u05ea\u05d9 \u05d2\u05d9\u05e9\u05d5\u05ea \u05dc\u05d4\u05ea\u05d2\u05d1\u05e8 \u05e2\u05dc \u05d6\u05d4:\",\n  \"\u05de\u05ea\u05d7\u05d9\u05dc\u05d9\u05dd \u05de\u05d0\u05d5\u05ea\u05d5 \u05d4\u05de\u05d5\u05d3\u05dc (\u05e9\u05e4\u05d4)\",\n  \"\u05dc\u05d0\u05de\u05df \u05de\u05d5\u05d3\u05dc\u05d9 reward \u05d6\u05d4\u05d9\u05dd \u05e2\u05dd \u05e8\u05d0\u05e9\u05d9\u05dd \u05dc\u05d9\u05e0\u05d0\u05e8\u05d9\u05d9\u05dd (\u05de\u05d0\u05d5\u05de\u05e0\u05d9\u05dd) \u05e9\u05d5\u05e0\u05d9\u05dd. \u05db\u05da \u05e6\u05e8\u05d9\u05da \u05dc\u05e9\u05de\u05d5\u05e8 \u05e8\u05e7 \u05de\u05d5\u05d3\u05dc \u05d0\u05d7\u05d3 \u05d5\u05d4\u05de\u05e9\u05e7\u05dc\u05d9\u05dd \u05e2\u05d1\u05d5\u05e8 \u05d4\u05e9\u05db\u05d1\u05d4 \u05d4\u05dc\u05d9\u05e0\u05d0\u05e8\u05d9\u05ea \u05e2\u05d1\u05d5\u05e8 \u05db\u05dc \u05de\u05d5\u05d3\u05dc.\",\n  \"\u05dc\u05d0\u05de\u05df \u05db\u05de\u05d4 \u05de\u05d5\u05d3\u05dc\u05d9 reward \u05d1\u05e9\u05d9\u05d8\u05d4 \u05e9\u05dc LoRa - \u05db\u05da \u05e0\u05e6\u05d8\u05e8\u05da \u05dc\u05e9\u05de\u05d5\u05e8 \u05e8\u05e7 \u05d0\u05ea \u05ea\u05d5\u05e1\u05e4\u05ea \u05d4\u05de\u05e9\u05e7\u05dc\u05d9\u05dd \u05dc\u05db\u05dc \u05e9\u05db\u05d1\u05d4 \u05e9\u05d6\u05d4 \u05d9\u05db\u05d5\u05dc \u05dc\u05d4\u05d9\u05d5\u05ea \u05d3\u05d9 \u05d6\u05d5\u05dc \u05de\u05d1\u05d7\u05d9\u05e0\u05ea \u05d4\u05de\u05e9\u05d0\u05d1\u05d9\u05dd\",\n  \"\u05d5\u05d0\u05d6 \u05d0\u05e4\u05e9\u05e8 \u05dc\u05e7\u05d7\u05ea \u05de\u05de\u05d5\u05e6\u05e2 \u05e9\u05dc \u05d4-rewards \u05e9\u05dc \u05db\u05dc \u05d4\u05de\u05d5\u05d3\u05dc\u05d9\u05dd \u05d0\u05d5 \u05d0\u05ea \u05d4\u05de\u05d9\u05e0\u05d5\u05de\u05d5\u05dd \u05d1\u05d9\u05e0\u05d9\u05d4\u05dd- \u05d9\u05e9 \u05dc\u05d0 \u05de\u05e2\u05d8 \u05d0\u05d5\u05e4\u05e6\u05d9\u05d5\u05ea\u2026\",\n  \"https://arxiv.org/abs/2401.16635\"\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nif (paragraphs.items.length !== 7) {\n  throw new Error(`Expected 7 paragraphs before edit, found ${paragraphs.items.length}`);\n}\n\n// Paragraphs 0-5: straightforward text replacement.\nfor (let i = 0; i <= 5; i++) {\n  paragraphs.items[i].insertText(newTexts[i], \"Replace\");\n}\n\n// The original last paragraph (index 6) held the arXiv link. Insert the\n// two new paragraphs right before it (in document order), then replace\n// its own text with the new link.\nconst lastParagraph = paragraphs.items[6];\nlastParagraph.insertParagraph(newTexts[6], \"Before\");\nlastParagraph.insertParagraph(newTexts[7], \"Before\");\nlastParagraph.insertText(newTexts[8], \"Replace\");\n\nawait context.sync();\n", "ps1": "# Apply the text edits described by the diff.\n# Paragraphs 1-6 (1-based COM indexing) get their text replaced in place;\n# two brand-new paragraphs are inserted right before the final (URL)\n# paragraph, whose own text is then updated to the new arXiv link.\n\n$newTexts = @(\n    '\u26a1\ufe0f\ud83d\ude80\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05d5\u05de\u05d9 \u05e9\u05dc \u05de\u05d9\u05d9\u05e7 25.06.24:\u26a1\ufe0f\ud83d\ude80',\n    'Improving Reinforcement Learning from Human Feedback with Efficient Reward Model Ensemble',\n    '\u05d4\u05e1\u05e7\u05d9\u05e8\u05d4 \u05d4\u05d6\u05d5 \u05de\u05de\u05e9\u05d9\u05db\u05d4 \u05d0\u05ea \u05e7\u05d5 \u05d4\u05e1\u05e7\u05d9\u05e8\u05d5\u05ea \u05e2\u05dc \u05d4\u05de\u05d0\u05de\u05e8\u05d9\u05dd \u05e9\u05de\u05e0\u05e1\u05d9\u05dd \u05dc\u05e9\u05e4\u05e8 \u05e9\u05d9\u05d8\u05d5\u05ea RLHF \u05dc\u05d8\u05d9\u05d5\u05d1 (instruction tuning \u05d0\u05d5 \u05e4\u05e9\u05d5\u05d8 fine-tuning) \u05e9\u05dc \u05de\u05d5\u05d3\u05dc\u05d9 \u05e9\u05e4\u05d4. \u05d1\u05d7\u05dc\u05e7 \u05e9\u05dc \u05e9\u05d9\u05d8\u05ea RLHF (\u05dc\u05de\u05e9\u05dc PPO) \u05d0\u05e0\u05d5 \u05de\u05d0\u05de\u05e0\u05d9\u05dd \u05de\u05d5\u05d3\u05dc reward \u05de\u05d1\u05d5\u05e1\u05e1 \u05e2\u05dc \u05e1\u05d8 \u05e9\u05dc \u05e9\u05d0\u05dc\u05d5\u05ea \u05d5\u05ea\u05e9\u05d5\u05d1\u05d5\u05ea \u05de\u05d3\u05d5\u05e8\u05d2\u05d5\u05ea \u05e2\u05dc \u05d9\u05d3\u05d9 \u05d4\u05de\u05ea\u05d9\u05d9\u05d2\u05d9\u05dd \u05d4\u05d0\u05e0\u05d5\u05e9\u05d9\u05d9\u05dd. \u05de\u05d8\u05e8\u05d4 \u05e9\u05dc \u05de\u05d5\u05d3\u05dc \u05d6\u05d4 \u05dc\u05e1\u05e4\u05e7 \u05e6\u05d9\u05d5\u05df \u05dc\u05d6\u05d5\u05d2 (\u05e9\u05d0\u05dc\u05d4, \u05ea\u05e9\u05d5\u05d1\u05d4) \u05db\u05d0\u05e9\u05e8 \u05e6\u05d9\u05d5\u05df \u05d2\u05d1\u05d5\u05d4 \u05de\u05e6\u05d1\u05d9\u05e2 \u05e2\u05dc \u05ea\u05e9\u05d5\u05d1\u05d4 \u05d8\u05d5\u05d1\u05d4 \u05d5\u05e8\u05e6\u05d5\u05d9\u05d4. \u05dc\u05d0\u05d7\u05e8 \u05db\u05df \u05d0\u05e0\u05d5 \u05de\u05d0\u05de\u05e0\u05d9\u05dd (\u05de\u05d8\u05d9\u05d9\u05d1\u05d9\u05dd) \u05de\u05d5\u05d3\u05dc \u05e9\u05e4\u05d4 \u05db\u05d0\u05e9\u05e8 \u05d4\u05de\u05d8\u05e8\u05d4 \u05d4\u05d9\u05d0 \u05de\u05e7\u05e1\u05d5\u05dd \u05e9\u05dc \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d4 reward \u05ea\u05d5\u05da \u05e9\u05de\u05d9\u05e8\u05ea \u05e9\u05dc \u05de\u05e9\u05e7\u05dc\u05d9 \u05d4\u05de\u05d5\u05d3\u05dc \u05dc\u05de\u05e9\u05e7\u05dc\u05d9\u05dd \u05e9\u05d4\u05ea\u05d7\u05dc\u05e0\u05d5 \u05de\u05d4\u05dd (\u05e0\u05de\u05d3\u05d3 \u05e2\u05dc \u05d9\u05d3\u05d9 KL divergence \u05d1\u05d9\u05df \u05d4\u05ea\u05e4\u05dc\u05d2\u05d5\u05d9\u05d5\u05ea \u05d4\u05d8\u05d5\u05e7\u05e0\u05d9\u05dd \u05e9\u05dc \u05e9\u05e0\u05d9 \u05d4\u05de\u05d5\u05d3\u05dc\u05d9\u05dd). \u05db\u05dc \u05d6\u05d4 \u05de\u05ea\u05d1\u05e6\u05e2 on-the-fly \u05db\u05d0\u05e9\u05e8 \u05d4\u05d3\u05d5\u05d2\u05de\u05d0\u05d5\u05ea \u05e0\u05d5\u05e6\u05e8\u05d5\u05ea \u05e2\u05dc\u05d9 \u05d9\u05d3\u05d9 \u05d4\u05d2\u05e8\u05e1\u05d4 \u05d4\u05e2\u05d3\u05db\u05e0\u05d9\u05ea \u05e9\u05dc \u05d4\u05de\u05d5\u05d3\u05dc \u05d1\u05de\u05d4\u05dc\u05da \u05d4\u05d0\u05d9\u05de\u05d5\u05df.',\n    '\u05d4\u05d1\u05e2\u05d9\u05d4 \u05e2\u05dd \u05d4\u05d2\u05d9\u05e9\u05d4 \u05d4\u05d9\u05d0 reward hacking \u05db\u05d0\u05e9\u05e8 \u05dc\u05de\u05e8\u05d5\u05ea \u05d0\u05d9\u05d1\u05e8 \u05d4\u05e8\u05d2\u05d5\u05dc\u05e8\u05d9\u05d6\u05e6\u05d9\u05d4 (KL) \u05d4\u05de\u05d5\u05d3\u05dc \u05de\u05ea\u05db\u05e0\u05e1 \u05dc\u05de\u05e9\u05e7\u05dc\u05d9\u05dd \u05e9\u05de\u05d2\u05d9\u05e2\u05d9\u05dd \u05dc\u05e2\u05e8\u05db\u05d9\u05dd \u05d2\u05d1\u05d5\u05d4\u05d9\u05dd \u05e9\u05dc \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d9\u05ea reward \u05db\u05d0\u05e9\u05e8 \u05d4\u05de\u05d5\u05d3\u05dc \u05e2\u05e6\u05de\u05d5 \u05f4\u05dc\u05d0 \u05de\u05e1\u05e4\u05e7 \u05d0\u05ea \u05d4\u05e1\u05d7\u05d5\u05e8\u05d4\u05f4. \u05d4\u05de\u05d0\u05de\u05e8 \u05de\u05e6\u05d9\u05e2 \u05dc\u05d4\u05e9\u05ea\u05de\u05e9 \u05d1\u05db\u05de\u05d4 \u05de\u05d5\u05d3\u05dc\u05d9 reward \u05db\u05d9 ensemble \u05d6\u05d4 \u05ea\u05de\u05d9\u05d3 \u05d8\u05d5\u05d1. \u05d4\u05d1\u05e2\u05d9\u05d4 \u05e9\u05dc\u05d4\u05d7\u05d6\u05d9\u05e7 \u05d9\u05d5\u05ea\u05e8 \u05de\u05de\u05d5\u05d3\u05dc \u05d0\u05d7\u05d3 \u05d1\u05d6\u05de\u05df \u05d4\u05d0\u05d9\u05de\u05d5\u05df \u05d6\u05d4 \u05d9\u05e7\u05e8 \u05de\u05d1\u05d7\u05d9\u05e0\u05ea \u05d4\u05de\u05e9\u05d0\u05d1\u05d9\u05dd. \u05d4\u05de\u05d0\u05de\u05e8 \u05de\u05e6\u05d9\u05e2 \u05e9\u05ea\u05d9 \u05d2\u05d9\u05e9\u05d5\u05ea \u05dc\u05d4\u05ea\u05d2\u05d1\u05e8 \u05e2\u05dc \u05d6\u05d4:',\n    '\u05de\u05ea\u05d7\u05d9\u05dc\u05d9\u05dd \u05de\u05d0\u05d5\u05ea\u05d5 \u05d4\u05de\u05d5\u05d3\u05dc (\u05e9\u05e4\u05d4)',\n    '\u05dc\u05d0\u05de\u05df \u05de\u05d5\u05d3\u05dc\u05d9 reward \u05d6\u05d4\u05d9\u05dd \u05e2\u05dd \u05e8\u05d0\u05e9\u05d9\u05dd \u05dc\u05d9\u05e0\u05d0\u05e8\u05d9\u05d9\u05dd (\u05de\u05d0\u05d5\u05de\u05e0\u05d9\u05dd) \u05e9\u05d5\u05e0\u05d9\u05dd. \u05db\u05da \u05e6\u05e8\u05d9\u05da \u05dc\u05e9\u05de\u05d5\u05e8 \u05e8\u05e7 \u05de\u05d5\u05d3\u05dc \u05d0\u05d7\u05d3 \u05d5\u05d4\u05de\u05e9\u05e7\u05dc\u05d9\u05dd \u05e2\u05d1\u05d5\u05e8 \u05d4\u05e9\u05db\u05d1\u05d4 \u05d4\u05dc\u05d9\u05e0\u05d0\u05e8\u05d9\u05ea \u05e2\u05d1\u05d5\u05e8 \u05db\u05dc \u05de\u05d5\u05d3\u05dc.',\n    '\u05dc\u05d0\u05de\u05df \u05db\u05de\u05d4 \u05de\u05d5\u05d3\u05dc\u05d9 reward \u05d1\u05e9\u05d9\u05d8\u05d4 \u05e9\u05dc LoRa - \u05db\u05da \u05e0\u05e6\u05d8\u05e8\u05da \u05dc\u05e9\u05de\u05d5\u05e8 \u05e8\u05e7 \u05d0\u05ea \u05ea\u05d5\u05e1\u05e4\u05ea \u05d4\u05de\u05e9\u05e7\u05dc\u05d9\u05dd \u05dc\u05db\u05dc \u05e9\u05db\u05d1\u05d4 \u05e9\u05d6\u05d4 \u05d9\u05db\u05d5\u05dc \u05dc\u05d4\u05d9\u05d5\u05ea \u05d3\u05d9 \u05d6\u05d5\u05dc \u05de\u05d1\u05d7\u05d9\u05e0\u05ea \u05d4\u05de\u05e9\u05d0\u05d1\u05d9\u05dd',\n    '\u05d5\u05d0\u05d6 \u05d0\u05e4\u05e9\u05e8 \u05dc\u05e7\u05d7\u05ea \u05de\u05de\u05d5\u05e6\u05e2 \u05e9\u05dc \u05d4-rewards \u05e9\u05dc \u05db\u05dc \u05d4\u05de\u05d5\u05d3\u05dc\u05d9\u05dd \u05d0\u05d5 \u05d0\u05ea \u05d4\u05de\u05d9\u05e0\u05d5\u05de\u05d5\u05dd \u05d1\u05d9\u05e0\u05d9\u05d4\u05dd- \u05d9\u05e9 \u05dc\u05d0 \u05de\u05e2\u05d8 \u05d0\u05d5\u05e4\u05e6\u05d9\u05d5\u05ea\u2026',\n    'https://arxiv.org/abs/2401.16635'\n)\n\n$d = $word.ActiveDocument\n\nif ($d.Paragraphs.Count -ne 7) {\n    throw \"Expected 7 paragraphs before edit, found $($d.Paragraphs.Count)\"\n}\n\n# Paragraphs 1-6 (1-based): straightforward text replacement.\nfor ($i = 1; $i -le 6; $i++) {\n    $d.Paragraphs($i).Range.Text = $newTexts[$i - 1]\n}\n\n# The original last paragraph (index 7) held the arXiv link. Insert the\n# two new paragraphs right before it (in document order), then replace\n# its own text with the new link.\n$lastRange = $d.Paragraphs($d.Paragraphs.Count).Range\n$lastRange.InsertParagraphBefore()\n$d.Paragraphs($d.Paragraphs.Count - 1).Range.Text = $newTexts[6]\n\n$lastRange = $d.Paragraphs($d.Paragraphs.Count).Range\n$lastRange.InsertParagraphBefore()\n$d.Paragraphs($d.Paragraphs.Count - 1).Range.Text = $newTexts[7]\n\n$d.Paragraphs($d.Paragraphs.Count).Range.Text = $newTexts[8]\n\n"}
